$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

# --- Slide Layout names (p:cSld/@name) ---
$layouts.Item(1).Name  = "Başlık Slaydı"            # slideLayout1.xml - Title Slide
$layouts.Item(2).Name  = "Başlık ve İçerik"          # slideLayout2.xml - Title and Content
$layouts.Item(3).Name  = "Bölüm Üstbilgisi"          # slideLayout3.xml - Section Header
$layouts.Item(4).Name  = "İki İçerik"                # slideLayout4.xml - Two Content
$layouts.Item(5).Name  = "Karşılaştırma"             # slideLayout5.xml - Comparison
$layouts.Item(6).Name  = "Yalnızca Başlık"           # slideLayout6.xml - Title Only
$layouts.Item(7).Name  = "Boş"                       # slideLayout7.xml - Blank
$layouts.Item(8).Name  = "Başlıklı İçerik"           # slideLayout8.xml - Content with Caption
$layouts.Item(9).Name  = "Başlıklı Resim"            # slideLayout9.xml - Picture with Caption
$layouts.Item(10).Name = "Başlık, Dikey Metin"       # slideLayout10.xml - Title and Vertical Text
$layouts.Item(11).Name = "Dikey Başlık ve Metin"     # slideLayout11.xml - Vertical Title and Text

# --- Shape rename on "Two Content" layout (slideLayout4.xml) ---
$layout4 = $layouts.Item(4)
try { $layout4.Shapes.Item(1).Name = "Başlık 1" } catch { }

# --- "Picture with Caption" layout (slideLayout9.xml) ---
$layout9 = $layouts.Item(9)
# Placeholder prompt text
$layout9.Shapes.Item(2).TextFrame.TextRange.Text = "Resim eklemek için simgeye tıklayın"
# Date placeholder shape rename
try { $layout9.Shapes.Item(4).Name = "Veri Yer Tutucusu 4" } catch { }
